$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("week3")

# C4: record 45 minutes spent (new value), formatted as time (h:mm) using the
# same black font already used elsewhere in the sheet (e.g. D18) so the
# workbook reuses the existing font instead of adding a new one.
$ws.Range("C4").Value = 45/1440
$ws.Range("C4").Font.Color = $ws.Range("D18").Font.Color
$ws.Range("C4").NumberFormat = "h:mm"

# E4: the stray note is no longer needed - clear it out entirely.
$ws.Range("E4").ClearContents()

# C7: record 30 minutes spent (previously blank).
$ws.Range("C7").Value = 30/1440

# Update the on-screen selection to match where the user ended up.
$ws.Range("C8").Select()
